$wb = $excel.ActiveWorkbook

# --------------------------------------------------------------------------
# 1) Update the "总计" (summary) sheet: shift the quarterly rows down by one
#    and insert a new "2022-Q3" row at the top of the data block.
# --------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)

# Shift columns B (label), C (count), D (value) down by one row, working
# from the bottom up so we never clobber a value before it is read.
for ($r = 8; $r -ge 2; $r--) {
    $wsTotal.Cells.Item($r + 1, 2).Value = $wsTotal.Cells.Item($r, 2).Value2
    $wsTotal.Cells.Item($r + 1, 3).Value = $wsTotal.Cells.Item($r, 3).Value2
    $wsTotal.Cells.Item($r + 1, 4).Value = $wsTotal.Cells.Item($r, 4).Value2
}

# New row 9 needs the same look as the rest of column A (bold/centered/bordered).
$wsTotal.Cells.Item(8, 1).Copy()
$wsTotal.Cells.Item(9, 1).PasteSpecial(-4122)
$wsTotal.Cells.Item(9, 1).Value = 7

# Fill the freed-up row 2 with the brand new 2022-Q3 summary figures.
$wsTotal.Cells.Item(2, 2).Value = "2022-Q3"
$wsTotal.Cells.Item(2, 3).Value = 12
$wsTotal.Cells.Item(2, 4).Value = 4.01

# --------------------------------------------------------------------------
# 2) Insert a brand new "2022-Q3" worksheet right after "总计" holding the
#    per-fund breakdown for the new quarter.
# --------------------------------------------------------------------------
$wsQ3 = $wb.Worksheets.Add($null, $wsTotal)
$wsQ3.Name = "2022-Q3"

# Match the page margins used throughout the rest of the workbook (the
# defaults for a brand new sheet differ slightly).
$wsQ3.PageSetup.LeftMargin = 54
$wsQ3.PageSetup.RightMargin = 54
$wsQ3.PageSetup.TopMargin = 72
$wsQ3.PageSetup.BottomMargin = 72
$wsQ3.PageSetup.HeaderMargin = 36
$wsQ3.PageSetup.FooterMargin = 36

# Make sure the numeric-looking identifiers / figures in columns B-G are
# stored as text (as they are on every other quarter sheet) so that values
# such as "012082" keep their leading zero.
$wsQ3.Range("B2:G13").NumberFormat = "@"

# Reuse the exact header style (bold, centered, bordered) already defined
# on the "总计" sheet for B1, and the column-A style used for the row index.
$wsTotal.Cells.Item(1, 2).Copy()
$wsQ3.Range("B1:H1").PasteSpecial(-4122)
$wsTotal.Cells.Item(2, 1).Copy()
$wsQ3.Range("A2:A13").PasteSpecial(-4122)

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $wsQ3.Cells.Item(1, 2 + $i).Value = $headers[$i]
}

$rows = @(
    @("519918", "华夏兴和混合", "45.93", "91.61", "5.07", "2.3287", 6),
    @("160311", "华夏蓝筹混合（LOF）A", "26.66", "91.52", "4.00", "1.0664", 9),
    @("012082", "博时数字经济18个月封闭混合A", "5.19", "96.61", "3.63", "0.1884", 6),
    @("001305", "九泰天富改革新动力混合A", "1.86", "94.71", "7.45", "0.1386", 7),
    @("217001", "招商安泰混合", "4.18", "75.08", "2.18", "0.0911", 10),
    @("001844", "九泰久益灵活配置混合C", "0.98", "93.32", "8.35", "0.0818", 4),
    @("001782", "九泰久益灵活配置混合A", "0.53", "93.32", "8.35", "0.0443", 4),
    @("014600", "博时回报严选混合A", "0.92", "92.53", "4.24", "0.0390", 6),
    @("009912", "九泰天富改革新动力混合C", "0.17", "94.71", "7.45", "0.0127", 7),
    @("012083", "博时数字经济18个月封闭混合C", "0.32", "96.61", "3.63", "0.0116", 6),
    @("014601", "博时回报严选混合C", "0.04", "92.53", "4.24", "0.0017", 6),
    @("015950", "华夏蓝筹混合（LOF）C", "0.04", "91.52", "4.00", "0.0016", 9)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $wsQ3.Cells.Item($r, 1).Value = $i
    $wsQ3.Cells.Item($r, 2).Value = $row[0]
    $wsQ3.Cells.Item($r, 3).Value = $row[1]
    $wsQ3.Cells.Item($r, 4).Value = $row[2]
    $wsQ3.Cells.Item($r, 5).Value = $row[3]
    $wsQ3.Cells.Item($r, 6).Value = $row[4]
    $wsQ3.Cells.Item($r, 7).Value = $row[5]
    $wsQ3.Cells.Item($r, 8).Value = $row[6]
}

# Leave the workbook's active sheet where it originally was ("总计").
$wsTotal.Activate()
